# Updates the crypto price/volume table with the latest scraped values.
# Price (column D) and Volume/1h change (column E) cells are stored as
# plain text (e.g. "1.705.39", "  -1.05%  "), so force a text number
# format before writing them to stop Excel from reinterpreting them as
# numbers/dates and mangling values such as "1.004" or "0.2660".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.299.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.705.39"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.05%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.94%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5308"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.14%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2660"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.87%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.44%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.73"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07633"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.76%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.496"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.10%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.709.25"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.78%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.940.51"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.02%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5778"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.73%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8146"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.60%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.61%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.306.91"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.77%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.25"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.84%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.610"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.77%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.24%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.963"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.12%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.27"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.82%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.698"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.73%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1200"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.67%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.207"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.87%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.10"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.75%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05370"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.73%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.286"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.48%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.469"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.92%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.402"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.83%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.30%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.867"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.86%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.415"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.51%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9460"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.57%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5806"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.24%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.31%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.779"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.85%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.004"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.08%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.041.42"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.84%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8404"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.12%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.08"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.49%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.848.76"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.94%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.85%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.73"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.27%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4518"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.97%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.11%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.036"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.41%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05227"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.94%  "
